# Update column G ("K") values for rows 2-13 on Sheet1, per the regenerated
# save_data (K instead of Strike#, std/mean, s_vals recalculation).

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$newValues = @{
    2  = 4
    3  = 2
    4  = 2
    5  = 5
    6  = 2
    7  = 1
    8  = 0
    9  = 1
    10 = 4
    11 = 6
    12 = 3
    13 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
